$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-13 02:00:20"
$wsZhCn.Range("G2").Value = "2016-01-13 02:01:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-13 02:00:45"
$wsDeDe.Range("G2").Value = "2016-01-13 02:02:20"
